$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new sensor readings were captured and prepended to the "falling"
# series, ahead of what is now row 4 (originally timestamp 0), so the rows
# that already existed slide down by two. The data window stays the same
# length, so the three oldest trailing readings (originally at rows 20-22,
# timestamps 1800/1900/2000) drop off the end.
$ws.Rows("2:3").Insert()

# Populate the two freshly inserted rows with their captured values.
$newRows = @(
    @(0,   "falling", -0.2734694480895995, 0.2277572751045226, -0.1111783366650344, 0.0204639863222837, -0.0009162978967650999, 0.007177666760981),
    @(100, "falling", -0.2649335861206055, 0.1057968139648434, -0.4681921228766454, 0.0155770638957619, -0.046578474342823, 0.0195476878434419)
)

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Drop the three oldest trailing samples that no longer fit in the window
# (these are the old rows 20-22, shifted down to 22-24 by the insert above).
$ws.Rows("22:24").Delete()

# Re-number the timestamp column so it stays a contiguous 100-unit cadence
# from 0 through to the last remaining row.
$ts = 200
for ($row = 4; $row -le 21; $row++) {
    $ws.Cells.Item($row, 1).Value = $ts
    $ts = $ts + 100
}
